# Re-sorts the weekly price rows (rows 2-20) for this chirimoya price
# subset. Only the per-record fields (Fecha, Calidad, Volumen, Precios,
# Unidad de comercialización, Origen, Precio $/Kg, Kg / unidad) move
# between rows; the descriptive columns (Mercado ID, Mercado, Región,
# Codreg, Tipo, Producto*, Categoría*, Variedad) are identical on every
# row, so only the varying columns need to be re-seated per the mapping
# below (new row -> source row, taken from the data's original order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new row number -> old row number that its data should come from
$rowMap = [ordered]@{
    2  = 4
    3  = 9
    4  = 10
    5  = 11
    6  = 7
    7  = 5
    8  = 6
    9  = 18
    10 = 15
    11 = 12
    12 = 17
    13 = 19
    14 = 20
    15 = 2
    16 = 14
    17 = 3
    18 = 8
    19 = 16
    20 = 13
}

$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot the current (pre-edit) values of the varying columns for every
# data row before any writes happen, so source rows aren't clobbered
# mid-loop by the permutation. Value2 is used (rather than Value) so dates
# come back as raw serials instead of formatted/variant wrappers.
$snapshot = @{}
for ($r = 2; $r -le 20; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Write each new row's values from its mapped source row's snapshot.
foreach ($newRow in $rowMap.Keys) {
    $oldRow = $rowMap[$newRow]
    $src = $snapshot[$oldRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value2 = $src[$col]
    }
}
